# Updated cryptos list cell values (price + 1h volume %) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.625.93"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "3.785.65"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.12"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.16"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").Value = "3.785.75"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").Value = "4.418.16"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "3.800.78"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.99"
$ws.Range("E17").Value = "  +5.97%  "
$ws.Range("D18").Value = "67.711.74"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.28"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.58"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.02"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.726"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000151"
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.67"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +1.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.32"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.91"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "3.934.30"
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.26"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.53"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.18"
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "3.747.23"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +7.55%  "
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.74"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.47"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "405.50"
$ws.Range("E48").Value = "  -4.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000279"
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.09"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +0.05%  "
